# Actualización automática del tracker — add two new result rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 14709122
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2025-09-19"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = "Juncheng Shang"
$ws.Range("D13").Value = "Brandon Nakashima"
$ws.Range("E13").Value = "Gana Brandon Nakashima"
$ws.Range("F13").Value = 1.62

# Row 14
$ws.Range("A14").Value = 14655133
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2025-09-18"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = "Daniel Masur"
$ws.Range("D14").Value = "Andres Martin"
$ws.Range("E14").Value = "Gana Daniel Masur"
$ws.Range("F14").Value = 2.63
